$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1054003333333333
$ws.Range("H2").Value = 0.316201
$ws.Range("I2").Value = 0.002842499753259756
$ws.Range("J2").Value = 0.002842499753259756
$ws.Range("M2").Value = 10.82167433333333
$ws.Range("N2").Value = 32.465023
$ws.Range("O2").Value = 0.09133543757015983
$ws.Range("P2").Value = 0.09133543757015983
$ws.Range("Q2").Value = 1.140608081958111
$ws.Range("R2").Value = 10.265472737623
$ws.Range("S2").Value = 0.0002596209587570511
$ws.Range("T2").Value = 0.0002596209587570511
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1054003333333333
$ws.Range("H3").Value = 0.316201
$ws.Range("I3").Value = 0.002842499753259756
$ws.Range("J3").Value = 0.002842499753259756
$ws.Range("M3").Value = 36.14140700000001
$ws.Range("O3").Value = 0.3050351656377608
$ws.Range("P3").Value = 0.3050351656377608
$ws.Range("Q3").Value = 3.809316344935668
$ws.Range("R3").Value = 34.28384710442101
$ws.Range("S3").Value = 0.0008670623830608838
$ws.Range("T3").Value = 0.0008670623830608838
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1054003333333333
$ws.Range("H4").Value = 0.316201
$ws.Range("I4").Value = 0.002842499753259756
$ws.Range("J4").Value = 0.002842499753259756
$ws.Range("M4").Value = 26.40107466666666
$ws.Range("N4").Value = 79.20322399999999
$ws.Range("O4").Value = 0.2228263051286729
$ws.Range("P4").Value = 0.2228263051286729
$ws.Range("Q4").Value = 2.782682070224888
$ws.Range("R4").Value = 25.044138632024
$ws.Range("S4").Value = 0.0006333837173480357
$ws.Range("T4").Value = 0.0006333837173480358
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.1054003333333333
$ws.Range("H5").Value = 0.316201
$ws.Range("I5").Value = 0.002842499753259756
$ws.Range("J5").Value = 0.002842499753259756
$ws.Range("M5").Value = 45.11859966666666
$ws.Range("N5").Value = 135.355799
$ws.Range("O5").Value = 0.3808030916634065
$ws.Range("P5").Value = 0.3808030916634065
$ws.Range("Q5").Value = 4.755515444399888
$ws.Range("R5").Value = 42.799638999599
$ws.Range("S5").Value = 0.001082432694093785
$ws.Range("T5").Value = 0.001082432694093785
$ws.Range("I6").Value = 0.5290853192840337
$ws.Range("J6").Value = 0.5290853192840336
$ws.Range("M6").Value = 10.82167433333333
$ws.Range("N6").Value = 32.465023
$ws.Range("O6").Value = 0.09133543757015983
$ws.Range("P6").Value = 0.09133543757015983
$ws.Range("Q6").Value = 212.305732139006
$ws.Range("R6").Value = 1910.751589251054
$ws.Range("S6").Value = 0.04832423914875494
$ws.Range("T6").Value = 0.04832423914875493
$ws.Range("I7").Value = 0.5290853192840337
$ws.Range("J7").Value = 0.5290853192840336
$ws.Range("M7").Value = 36.14140700000001
$ws.Range("O7").Value = 0.3050351656377608
$ws.Range("P7").Value = 0.3050351656377608
$ws.Range("Q7").Value = 709.0425785623622
$ws.Range("R7").Value = 6381.38320706126
$ws.Range("S7").Value = 0.1613896280043128
$ws.Range("T7").Value = 0.1613896280043127
$ws.Range("I8").Value = 0.5290853192840337
$ws.Range("J8").Value = 0.5290853192840336
$ws.Range("M8").Value = 26.40107466666666
$ws.Range("N8").Value = 79.20322399999999
$ws.Range("O8").Value = 0.2228263051286729
$ws.Range("P8").Value = 0.2228263051286729
$ws.Range("Q8").Value = 517.9512258189279
$ws.Range("R8").Value = 4661.561032370352
$ws.Range("S8").Value = 0.1178941267938854
$ws.Range("T8").Value = 0.1178941267938854
$ws.Range("I9").Value = 0.5290853192840337
$ws.Range("J9").Value = 0.5290853192840336
$ws.Range("M9").Value = 45.11859966666666
$ws.Range("N9").Value = 135.355799
$ws.Range("O9").Value = 0.3808030916634065
$ws.Range("P9").Value = 0.3808030916634065
$ws.Range("Q9").Value = 885.1622253880779
$ws.Range("R9").Value = 7966.460028492702
$ws.Range("S9").Value = 0.2014773253370806
$ws.Range("T9").Value = 0.2014773253370805
$ws.Range("G10").Value = 1.502894666666666
$ws.Range("H10").Value = 4.508684
$ws.Range("I10").Value = 0.04053096972345505
$ws.Range("J10").Value = 0.04053096972345505
$ws.Range("M10").Value = 10.82167433333333
$ws.Range("N10").Value = 32.465023
$ws.Range("O10").Value = 0.09133543757015983
$ws.Range("P10").Value = 0.09133543757015983
$ws.Range("Q10").Value = 16.26383663997022
$ws.Range("R10").Value = 146.374529759732
$ws.Range("S10").Value = 0.003701913854834666
$ws.Range("T10").Value = 0.003701913854834666
$ws.Range("G11").Value = 1.502894666666666
$ws.Range("H11").Value = 4.508684
$ws.Range("I11").Value = 0.04053096972345505
$ws.Range("J11").Value = 0.04053096972345505
$ws.Range("M11").Value = 36.14140700000001
$ws.Range("O11").Value = 0.3050351656377608
$ws.Range("P11").Value = 0.3050351656377608
$ws.Range("Q11").Value = 54.31672782612934
$ws.Range("R11").Value = 488.8505504351641
$ws.Range("S11").Value = 0.01236337106305318
$ws.Range("T11").Value = 0.01236337106305318
$ws.Range("G12").Value = 1.502894666666666
$ws.Range("H12").Value = 4.508684
$ws.Range("I12").Value = 0.04053096972345505
$ws.Range("J12").Value = 0.04053096972345505
$ws.Range("M12").Value = 26.40107466666666
$ws.Range("N12").Value = 79.20322399999999
$ws.Range("O12").Value = 0.2228263051286729
$ws.Range("P12").Value = 0.2228263051286729
$ws.Range("Q12").Value = 39.67803431080177
$ws.Range("R12").Value = 357.1023087972159
$ws.Range("S12").Value = 0.009031366226759596
$ws.Range("T12").Value = 0.009031366226759598
$ws.Range("G13").Value = 1.502894666666666
$ws.Range("H13").Value = 4.508684
$ws.Range("I13").Value = 0.04053096972345505
$ws.Range("J13").Value = 0.04053096972345505
$ws.Range("M13").Value = 45.11859966666666
$ws.Range("N13").Value = 135.355799
$ws.Range("O13").Value = 0.3808030916634065
$ws.Range("P13").Value = 0.3808030916634065
$ws.Range("Q13").Value = 67.80850280650176
$ws.Range("R13").Value = 610.2765252585159
$ws.Range("S13").Value = 0.0154343185788076
$ws.Range("T13").Value = 0.01543431857880761
$ws.Range("G14").Value = 15.85329466666666
$ws.Range("H14").Value = 47.559884
$ws.Range("I14").Value = 0.4275412112392517
$ws.Range("J14").Value = 0.4275412112392516
$ws.Range("M14").Value = 10.82167433333333
$ws.Range("N14").Value = 32.465023
$ws.Range("O14").Value = 0.09133543757015983
$ws.Range("P14").Value = 0.09133543757015983
$ws.Range("Q14").Value = 171.5591919930369
$ws.Range("R14").Value = 1544.032727937332
$ws.Range("S14").Value = 0.03904966360781319
$ws.Range("T14").Value = 0.03904966360781318
$ws.Range("G15").Value = 15.85329466666666
$ws.Range("H15").Value = 47.559884
$ws.Range("I15").Value = 0.4275412112392517
$ws.Range("J15").Value = 0.4275412112392516
$ws.Range("M15").Value = 36.14140700000001
$ws.Range("O15").Value = 0.3050351656377608
$ws.Range("P15").Value = 0.3050351656377608
$ws.Range("Q15").Value = 572.9603748389294
$ws.Range("R15").Value = 5156.643373550364
$ws.Range("S15").Value = 0.130415104187334
$ws.Range("T15").Value = 0.130415104187334
$ws.Range("G16").Value = 15.85329466666666
$ws.Range("H16").Value = 47.559884
$ws.Range("I16").Value = 0.4275412112392517
$ws.Range("J16").Value = 0.4275412112392516
$ws.Range("M16").Value = 26.40107466666666
$ws.Range("N16").Value = 79.20322399999999
$ws.Range("O16").Value = 0.2228263051286729
$ws.Range("P16").Value = 0.2228263051286729
$ws.Range("Q16").Value = 418.544016207335
$ws.Range("R16").Value = 3766.896145866016
$ws.Range("S16").Value = 0.09526742839067988
$ws.Range("T16").Value = 0.09526742839067988
$ws.Range("G17").Value = 15.85329466666666
$ws.Range("H17").Value = 47.559884
$ws.Range("I17").Value = 0.4275412112392517
$ws.Range("J17").Value = 0.4275412112392516
$ws.Range("M17").Value = 45.11859966666666
$ws.Range("N17").Value = 135.355799
$ws.Range("O17").Value = 0.3808030916634065
$ws.Range("P17").Value = 0.3808030916634065
$ws.Range("Q17").Value = 715.2784554630349
$ws.Range("R17").Value = 6437.506099167315
$ws.Range("S17").Value = 0.1628090150534246
$ws.Range("T17").Value = 0.1628090150534246
